$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (shifts existing rows 36-44 down to 37-45),
# matching the diff where a new "FY_4.png" entry is inserted into the table
# and the table dimension grows from A1:G44 to A1:G45.
$ws.Rows.Item(36).Insert()

# Re-apply the full B:G value set for every data row (the underlying detection
# run was re-scored, so most numeric columns shift slightly) plus the new row's
# A/F/G values and the handful of F/G (name / correctness) calls that flipped.

$ws.Range("A1").Value = "AK_1.png"
$ws.Range("B1").Value = 0.351
$ws.Range("C1").Value = 0.001
$ws.Range("D1").Value = 0.451
$ws.Range("E1").Value = 0.467
$ws.Range("F1").Value = "Tidak Diketahui"
$ws.Range("G1").Value = "Salah"

$ws.Range("A2").Value = "AK_2.png"
$ws.Range("B2").Value = 0.6860000000000001
$ws.Range("C2").Value = 0.002
$ws.Range("D2").Value = 0.398
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "Akhlak Kamiswara"
$ws.Range("G2").Value = "Benar"

$ws.Range("A3").Value = "AK_3.png"
$ws.Range("B3").Value = 0.677
$ws.Range("C3").Value = 0.002
$ws.Range("D3").Value = 0.528
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Akhlak Kamiswara"
$ws.Range("G3").Value = "Benar"

$ws.Range("A4").Value = "AK_4.png"
$ws.Range("B4").Value = 0.334
$ws.Range("C4").Value = 0.001
$ws.Range("D4").Value = 0.43
$ws.Range("E4").Value = 0.533
$ws.Range("F4").Value = "Akhlak Kamiswara"
$ws.Range("G4").Value = "Benar"

$ws.Range("A5").Value = "AK_5.png"
$ws.Range("B5").Value = 0.32
$ws.Range("C5").Value = 0.001
$ws.Range("D5").Value = 0.439
$ws.Range("E5").Value = 0.667
$ws.Range("F5").Value = "Akhlak Kamiswara"
$ws.Range("G5").Value = "Benar"

$ws.Range("A6").Value = "MIB_1.png"
$ws.Range("B6").Value = 1.343
$ws.Range("C6").Value = 0.004
$ws.Range("D6").Value = 0.08599999999999999
$ws.Range("E6").Value = 0.467
$ws.Range("F6").Value = "Tidak Diketahui"
$ws.Range("G6").Value = "Salah"

$ws.Range("A7").Value = "MIB_2.png"
$ws.Range("B7").Value = 1.03
$ws.Range("C7").Value = 0.003
$ws.Range("D7").Value = 0.311
$ws.Range("E7").Value = 0.867
$ws.Range("F7").Value = "Muhammad Iqbal Baqi"
$ws.Range("G7").Value = "Benar"

$ws.Range("A8").Value = "MIB_3.png"
$ws.Range("B8").Value = 1.568
$ws.Range("C8").Value = 0.005
$ws.Range("D8").Value = 0.627
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "Muhammad Iqbal Baqi"
$ws.Range("G8").Value = "Benar"

$ws.Range("A9").Value = "MIB_4.png"
$ws.Range("B9").Value = 0.987
$ws.Range("C9").Value = 0.003
$ws.Range("D9").Value = 0.223
$ws.Range("E9").Value = 0.9330000000000001
$ws.Range("F9").Value = "Muhammad Iqbal Baqi"
$ws.Range("G9").Value = "Benar"

$ws.Range("A10").Value = "MIB_5.png"
$ws.Range("B10").Value = 1.378
$ws.Range("C10").Value = 0.005
$ws.Range("D10").Value = 0.247
$ws.Range("E10").Value = 0.9330000000000001
$ws.Range("F10").Value = "Muhammad Iqbal Baqi"
$ws.Range("G10").Value = "Benar"

$ws.Range("A11").Value = "AAH_1.png"
$ws.Range("B11").Value = 0.6850000000000001
$ws.Range("C11").Value = 0.002
$ws.Range("D11").Value = 0.265
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Andrea Ayunove Hutami"
$ws.Range("G11").Value = "Benar"

$ws.Range("A12").Value = "AAH_2.png"
$ws.Range("B12").Value = 0.949
$ws.Range("C12").Value = 0.003
$ws.Range("D12").Value = 0.681
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "Andrea Ayunove Hutami"
$ws.Range("G12").Value = "Benar"

$ws.Range("A13").Value = "AAH_3.png"
$ws.Range("B13").Value = 0.718
$ws.Range("C13").Value = 0.002
$ws.Range("D13").Value = 0.281
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Andrea Ayunove Hutami"
$ws.Range("G13").Value = "Benar"

$ws.Range("A14").Value = "TI_1.png"
$ws.Range("B14").Value = 0.759
$ws.Range("C14").Value = 0.003
$ws.Range("D14").Value = 0.488
$ws.Range("E14").Value = 0.533
$ws.Range("F14").Value = "Toni Ismail"
$ws.Range("G14").Value = "Benar"

$ws.Range("A15").Value = "TI_2.png"
$ws.Range("B15").Value = 0.791
$ws.Range("C15").Value = 0.003
$ws.Range("D15").Value = 0.413
$ws.Range("E15").Value = 0.667
$ws.Range("F15").Value = "Toni Ismail"
$ws.Range("G15").Value = "Benar"

$ws.Range("A16").Value = "TI_3.png"
$ws.Range("B16").Value = 0.5590000000000001
$ws.Range("C16").Value = 0.002
$ws.Range("D16").Value = 0.676
$ws.Range("E16").Value = 0.667
$ws.Range("F16").Value = "Toni Ismail"
$ws.Range("G16").Value = "Benar"

$ws.Range("A17").Value = "TI_4.png"
$ws.Range("B17").Value = 0.532
$ws.Range("C17").Value = 0.002
$ws.Range("D17").Value = 0.383
$ws.Range("E17").Value = 0.8
$ws.Range("F17").Value = "Toni Ismail"
$ws.Range("G17").Value = "Benar"

$ws.Range("A18").Value = "TI_5.png"
$ws.Range("B18").Value = 0.88
$ws.Range("C18").Value = 0.003
$ws.Range("D18").Value = 0.449
$ws.Range("E18").Value = 0.667
$ws.Range("F18").Value = "Toni Ismail"
$ws.Range("G18").Value = "Benar"

$ws.Range("A19").Value = "RAS_1.png"
$ws.Range("B19").Value = 0.469
$ws.Range("C19").Value = 0.002
$ws.Range("D19").Value = 0.445
$ws.Range("E19").Value = 0.533
$ws.Range("F19").Value = "Ridha Ayu Salsabila"
$ws.Range("G19").Value = "Benar"

$ws.Range("A20").Value = "RAS_2.png"
$ws.Range("B20").Value = 0.864
$ws.Range("C20").Value = 0.003
$ws.Range("D20").Value = 0.294
$ws.Range("E20").Value = 0.867
$ws.Range("F20").Value = "Ridha Ayu Salsabila"
$ws.Range("G20").Value = "Benar"

$ws.Range("A21").Value = "RAS_3.png"
$ws.Range("B21").Value = 0.372
$ws.Range("C21").Value = 0.001
$ws.Range("D21").Value = 0.343
$ws.Range("E21").Value = 0.467
$ws.Range("F21").Value = "Tidak Diketahui"
$ws.Range("G21").Value = "Salah"

$ws.Range("A22").Value = "RAS_4.png"
$ws.Range("B22").Value = 1.075
$ws.Range("C22").Value = 0.004
$ws.Range("D22").Value = 0.169
$ws.Range("E22").Value = 0.4
$ws.Range("F22").Value = "Tidak Diketahui"
$ws.Range("G22").Value = "Salah"

$ws.Range("A23").Value = "RAS_5.png"
$ws.Range("B23").Value = 1.009
$ws.Range("C23").Value = 0.003
$ws.Range("D23").Value = 0.362
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "Ridha Ayu Salsabila"
$ws.Range("G23").Value = "Benar"

$ws.Range("A24").Value = "RR_1.png"
$ws.Range("B24").Value = 1.302
$ws.Range("C24").Value = 0.004
$ws.Range("D24").Value = 0.482
$ws.Range("E24").Value = 0.9330000000000001
$ws.Range("F24").Value = "Rafiqo Rapitasari"
$ws.Range("G24").Value = "Benar"

$ws.Range("A25").Value = "RR_2.png"
$ws.Range("B25").Value = 1.265
$ws.Range("C25").Value = 0.004
$ws.Range("D25").Value = 0.52
$ws.Range("E25").Value = 0.9330000000000001
$ws.Range("F25").Value = "Rafiqo Rapitasari"
$ws.Range("G25").Value = "Benar"

$ws.Range("A26").Value = "RR_3.png"
$ws.Range("B26").Value = 0.926
$ws.Range("C26").Value = 0.003
$ws.Range("D26").Value = 0.106
$ws.Range("E26").Value = 0.867
$ws.Range("F26").Value = "Rafiqo Rapitasari"
$ws.Range("G26").Value = "Benar"

$ws.Range("A27").Value = "RR_4.png"
$ws.Range("B27").Value = 1.228
$ws.Range("C27").Value = 0.004
$ws.Range("D27").Value = 0.526
$ws.Range("E27").Value = 0.867
$ws.Range("F27").Value = "Rafiqo Rapitasari"
$ws.Range("G27").Value = "Benar"

$ws.Range("A28").Value = "RR_5.png"
$ws.Range("B28").Value = 1.282
$ws.Range("C28").Value = 0.004
$ws.Range("D28").Value = 0.544
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = "Rafiqo Rapitasari"
$ws.Range("G28").Value = "Benar"

$ws.Range("A29").Value = "AR_1.png"
$ws.Range("B29").Value = 0.679
$ws.Range("C29").Value = 0.002
$ws.Range("D29").Value = 0.397
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Arizli Romadhon"
$ws.Range("G29").Value = "Benar"

$ws.Range("A30").Value = "GA_1.png"
$ws.Range("B30").Value = 1.439
$ws.Range("C30").Value = 0.005
$ws.Range("D30").Value = 0.55
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Gege Ardiyansyah"
$ws.Range("G30").Value = "Benar"

$ws.Range("A31").Value = "GA_2.png"
$ws.Range("B31").Value = 0.496
$ws.Range("C31").Value = 0.002
$ws.Range("D31").Value = 0.275
$ws.Range("E31").Value = 0.9330000000000001
$ws.Range("F31").Value = "Gege Ardiyansyah"
$ws.Range("G31").Value = "Benar"

$ws.Range("A32").Value = "GA_3.png"
$ws.Range("B32").Value = 0.5600000000000001
$ws.Range("C32").Value = 0.002
$ws.Range("D32").Value = 0.188
$ws.Range("E32").Value = 0.733
$ws.Range("F32").Value = "Gege Ardiyansyah"
$ws.Range("G32").Value = "Benar"

$ws.Range("A33").Value = "FY_1.png"
$ws.Range("B33").Value = 0.986
$ws.Range("C33").Value = 0.003
$ws.Range("D33").Value = 0.308
$ws.Range("E33").Value = 0.733
$ws.Range("F33").Value = "Fanny Yusuf"
$ws.Range("G33").Value = "Benar"

$ws.Range("A34").Value = "FY_2.png"
$ws.Range("B34").Value = 1.478
$ws.Range("C34").Value = 0.005
$ws.Range("D34").Value = 0.481
$ws.Range("E34").Value = 0.6
$ws.Range("F34").Value = "Fanny Yusuf"
$ws.Range("G34").Value = "Benar"

$ws.Range("A35").Value = "FY_3.png"
$ws.Range("B35").Value = 1.381
$ws.Range("C35").Value = 0.005
$ws.Range("D35").Value = 0.468
$ws.Range("E35").Value = 0.6
$ws.Range("F35").Value = "Fanny Yusuf"
$ws.Range("G35").Value = "Benar"

$ws.Range("A36").Value = "FY_4.png"
$ws.Range("B36").Value = 1.196
$ws.Range("C36").Value = 0.004
$ws.Range("D36").Value = 0.376
$ws.Range("E36").Value = 0.6
$ws.Range("F36").Value = "Fanny Yusuf"
$ws.Range("G36").Value = "Benar"

$ws.Range("A37").Value = "TO_1.png"
$ws.Range("B37").Value = 0.624
$ws.Range("C37").Value = 0.002
$ws.Range("D37").Value = 0.397
$ws.Range("E37").Value = 0.8
$ws.Range("F37").Value = "Tiara Oktavian"
$ws.Range("G37").Value = "Benar"

$ws.Range("A38").Value = "TO_2.png"
$ws.Range("B38").Value = 0.898
$ws.Range("C38").Value = 0.003
$ws.Range("D38").Value = 0.374
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = "Tiara Oktavian"
$ws.Range("G38").Value = "Benar"

$ws.Range("A39").Value = "TO_3.png"
$ws.Range("B39").Value = 0.657
$ws.Range("C39").Value = 0.002
$ws.Range("D39").Value = 0.391
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = "Tiara Oktavian"
$ws.Range("G39").Value = "Benar"

$ws.Range("A40").Value = "TO_4.png"
$ws.Range("B40").Value = 2.874
$ws.Range("C40").Value = 0.008999999999999999
$ws.Range("D40").Value = 0.121
$ws.Range("E40").Value = 0.8
$ws.Range("F40").Value = "Tiara Oktavian"
$ws.Range("G40").Value = "Benar"

$ws.Range("A41").Value = "TO_5.png"
$ws.Range("B41").Value = 2.803
$ws.Range("C41").Value = 0.008999999999999999
$ws.Range("D41").Value = 0.091
$ws.Range("E41").Value = 0.9330000000000001
$ws.Range("F41").Value = "Tiara Oktavian"
$ws.Range("G41").Value = "Benar"

$ws.Range("A42").Value = "TD_1.png"
$ws.Range("B42").Value = 2.132
$ws.Range("C42").Value = 0.007
$ws.Range("D42").Value = 0.033
$ws.Range("E42").Value = 0.333
$ws.Range("F42").Value = "Tidak Diketahui"
$ws.Range("G42").Value = "Benar"

$ws.Range("A43").Value = "TD_2.png"
$ws.Range("B43").Value = 2.297
$ws.Range("C43").Value = 0.008
$ws.Range("D43").Value = 0.059
$ws.Range("E43").Value = 0.267
$ws.Range("F43").Value = "Tidak Diketahui"
$ws.Range("G43").Value = "Benar"

$ws.Range("A44").Value = "TD_3.png"
$ws.Range("B44").Value = 0.91
$ws.Range("C44").Value = 0.003
$ws.Range("D44").Value = 0.193
$ws.Range("E44").Value = 0.333
$ws.Range("F44").Value = "Tidak Diketahui"
$ws.Range("G44").Value = "Benar"

$ws.Range("A45").Value = "TD_4.png"
$ws.Range("B45").Value = 0.796
$ws.Range("C45").Value = 0.003
$ws.Range("D45").Value = 0.105
$ws.Range("E45").Value = 0.267
$ws.Range("F45").Value = "Tidak Diketahui"
$ws.Range("G45").Value = "Benar"
